# ---------------------------------------------------------------------------
# Adds the "Right / Wrong / Not Attempt / Max" summary block (rows 8-14) and
# the per-question "Student Ans / Correct Ans" tables (rows 15-40) to the
# quiz mark-sheet, colouring each student answer green when it matches the
# correct answer and red when it doesn't.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteSpecial "formats only" paste type used throughout.
$xlPasteFormats = -4122
# xlCenter alignment constant.
$xlCenter = -4108

# ---------------------------------------------------------------------------
# Step 1: build five reusable "template" styles in a scratch area (column Z)
# by cloning the look of existing cells (Century font, correct size, no
# stray theme colour) and adding the thin border + centred alignment that
# every new cell in this block needs. Each template cell is later copied
# (format only) onto the real destination ranges, and the scratch area is
# cleared at the very end.
# ---------------------------------------------------------------------------

# Plain style (Century 12, no font colour) -> mirrors A6's look.
$ws.Range("A6").Copy()
$ws.Range("Z1").PasteSpecial($xlPasteFormats)
$ws.Range("Z1").Borders.LineStyle = 1
$ws.Range("Z1").HorizontalAlignment = $xlCenter

# Green style (Century 12, green font) -> correct-answer highlight.
$ws.Range("A6").Copy()
$ws.Range("Z2").PasteSpecial($xlPasteFormats)
$ws.Range("Z2").Borders.LineStyle = 1
$ws.Range("Z2").HorizontalAlignment = $xlCenter
$ws.Range("Z2").Font.Color = 32768

# Red style (Century 12, red font) -> wrong-answer / "Wrong" column highlight.
$ws.Range("A6").Copy()
$ws.Range("Z3").PasteSpecial($xlPasteFormats)
$ws.Range("Z3").Borders.LineStyle = 1
$ws.Range("Z3").HorizontalAlignment = $xlCenter
$ws.Range("Z3").Font.Color = 255

# Blue style (Century 12, blue font) -> "Correct Ans" column.
$ws.Range("A6").Copy()
$ws.Range("Z4").PasteSpecial($xlPasteFormats)
$ws.Range("Z4").Borders.LineStyle = 1
$ws.Range("Z4").HorizontalAlignment = $xlCenter
$ws.Range("Z4").Font.Color = 16711680

# Bold style (Century 12 bold, no font colour) -> mirrors B6's look, used
# for the "Student Ans" / "Correct Ans" column headers.
$ws.Range("B6").Copy()
$ws.Range("Z5").PasteSpecial($xlPasteFormats)
$ws.Range("Z5").Borders.LineStyle = 1
$ws.Range("Z5").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------------
# Step 2: blank spacer rows (8, 13, 14) - present but unstyled/empty.
# ---------------------------------------------------------------------------

$ws.Range("A8:E8").Value = ""
$ws.Range("A13:E13").Value = ""
$ws.Range("A14:E14").Value = ""

# ---------------------------------------------------------------------------
# Step 3: "Right / Wrong / Not Attempt / Max" summary table (rows 9-12).
# ---------------------------------------------------------------------------

# Row 9 (column headers) is entirely "plain" styled - including B9/C9.
$ws.Range("Z1").Copy()
$ws.Range("A9:E12").PasteSpecial($xlPasteFormats)

# Rows 10-12 colour the "Right"/"Wrong" number columns green/red.
$ws.Range("Z2").Copy()
$ws.Range("B10:B12").PasteSpecial($xlPasteFormats)
$ws.Range("Z3").Copy()
$ws.Range("C10:C12").PasteSpecial($xlPasteFormats)

$ws.Range("B9").Value = "Right"
$ws.Range("C9").Value = "Wrong"
$ws.Range("D9").Value = "Not Attempt"
$ws.Range("E9").Value = "Max"

$ws.Range("A10").Value = "No."
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 12
$ws.Range("D10").Value = 7
$ws.Range("E10").Value = 28

$ws.Range("A11").Value = "Marking"
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0

$ws.Range("A12").Value = "Total"
$ws.Range("B12").Value = 45
$ws.Range("C12").Value = -12

$ws.Range("Z4").Copy()
$ws.Range("E12").PasteSpecial($xlPasteFormats)
$ws.Range("E12").Value = "33/140"

# ---------------------------------------------------------------------------
# Step 4: "Student Ans" / "Correct Ans" headers (row 15) for both tables.
# ---------------------------------------------------------------------------

$ws.Range("Z5").Copy()
$ws.Range("A15:B15").PasteSpecial($xlPasteFormats)
$ws.Range("D15:E15").PasteSpecial($xlPasteFormats)

$ws.Range("A15").Value = "Student Ans"
$ws.Range("B15").Value = "Correct Ans"
$ws.Range("D15").Value = "Student Ans"
$ws.Range("E15").Value = "Correct Ans"

# ---------------------------------------------------------------------------
# Step 5: per-question answer tables. Column A/D = student answer (green if
# it matches the correct answer, red otherwise); column B/E = correct answer
# (always blue).
# ---------------------------------------------------------------------------

$table1 = @(
    @{row=16; a="Option A"; b="Option A"},
    @{row=17; a="Option D"; b="Option D"},
    @{row=18; a="Option B"; b="Option B"},
    @{row=19; a="Option C"; b="Option C"},
    @{row=20; a=$null; b="Option B"},
    @{row=21; a="Option D"; b="Option C"},
    @{row=22; a="Option D"; b="Option D"},
    @{row=23; a="Option C"; b="Option D"},
    @{row=24; a=$null; b="Option A"},
    @{row=25; a="Option D"; b="Option A"},
    @{row=26; a="Option B"; b="Option C"},
    @{row=27; a="Option D"; b="Option A"},
    @{row=28; a="Option C"; b="Option D"},
    @{row=29; a="Option C"; b="Option D"},
    @{row=30; a="Option B"; b="Option B"},
    @{row=31; a="Option A"; b="Option D"},
    @{row=32; a="Option B"; b="Option C"},
    @{row=33; a="Option D"; b="Option D"},
    @{row=34; a=$null; b="Option B"},
    @{row=35; a=$null; b="Option D"},
    @{row=36; a=$null; b="Option A"},
    @{row=37; a="Option B"; b="Option A"},
    @{row=38; a=$null; b="Option A"},
    @{row=39; a="Option C"; b="Option D"},
    @{row=40; a=$null; b="Option D"},
)

$table2 = @(
    @{row=16; a="Option A"; b="Option A"},
    @{row=17; a="Option B"; b="Option C"},
    @{row=18; a="Option D"; b="Option D"},
)

foreach ($d in $table1) {
    $cellA = $ws.Cells.Item($d.row, 1)
    if ($d.a -eq $d.b) {
        $ws.Range("Z2").Copy()
    } else {
        $ws.Range("Z3").Copy()
    }
    $cellA.PasteSpecial($xlPasteFormats)
    if ($d.a -ne $null) {
        $cellA.Value = $d.a
    } else {
        $cellA.Value = ""
    }

    $cellB = $ws.Cells.Item($d.row, 2)
    $ws.Range("Z4").Copy()
    $cellB.PasteSpecial($xlPasteFormats)
    $cellB.Value = $d.b
}

foreach ($d in $table2) {
    $cellD = $ws.Cells.Item($d.row, 4)
    if ($d.a -eq $d.b) {
        $ws.Range("Z2").Copy()
    } else {
        $ws.Range("Z3").Copy()
    }
    $cellD.PasteSpecial($xlPasteFormats)
    if ($d.a -ne $null) {
        $cellD.Value = $d.a
    } else {
        $cellD.Value = ""
    }

    $cellE = $ws.Cells.Item($d.row, 5)
    $ws.Range("Z4").Copy()
    $cellE.PasteSpecial($xlPasteFormats)
    $cellE.Value = $d.b
}

# ---------------------------------------------------------------------------
# Step 6: drop the scratch template cells so they don't leak into the
# worksheet's used range.
# ---------------------------------------------------------------------------

$ws.Range("Z1:Z5").Clear()
